$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 12:27"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4813984
$ws.Range("C4").Value = 337
$ws.Range("D4").Value = 2380548
$ws.Range("E4").Value = 2275064
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 158372

# Iran (row 14)
$ws.Range("B14").Value = 312035
$ws.Range("C14").Value = 2598
$ws.Range("D14").Value = 270228
$ws.Range("E14").Value = 24402
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 215
$ws.Range("H14").Value = 17405

# Banglades (row 19)
$ws.Range("B19").Value = 242102
$ws.Range("C19").Value = 1356
$ws.Range("D19").Value = 137905
$ws.Range("E19").Value = 101013
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 3184

# Countries list re-sorted: Kuwait now placed before Bielorrusia.
# Row 40 becomes Kuwait (new, updated figures); row 41 becomes Bielorrusia
# (figures unchanged from before the edit).
$ws.Range("A40").Value = "Kuwait"
$ws.Range("B40").Value = 68229
$ws.Range("C40").Value = 318
$ws.Range("D40").Value = 59739
$ws.Range("E40").Value = 8029
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 461

$ws.Range("A41").Value = "Bielorrusia"
$ws.Range("B41").Value = 68067
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 62896
$ws.Range("E41").Value = 4604
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 567

# Rumania (row 45)
$ws.Range("B45").Value = 54009
$ws.Range("C45").Value = 823
$ws.Range("D45").Value = 27750
$ws.Range("E45").Value = 23827
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 19
$ws.Range("H45").Value = 2432

# Malasia (row 88)
$ws.Range("B88").Value = 9001
$ws.Range("C88").Value = 2
$ws.Range("D88").Value = 8668
$ws.Range("E88").Value = 208
